# Implement revenue figures on the manager homepage: fill in the
# (previously blank) order-detail row for order #2 and add the order
# number for the next row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order #2 (row 2): name / email placeholders, pizza counts and total revenue.
$ws.Range("B2").Value = "er"
$ws.Range("C2").Value = "e"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 36.47

# Order #3 (row 3) now has an order number.
$ws.Range("A3").Value = 3
